$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Periodo Mora" column (E16:E30). The database of account
# statements was re-sorted into ascending period order (2109 .. 2211),
# while the row-level formatting (the thin-bordered "first row" /
# "last row" styles already baked into rows 16 and 30) stays put - only
# the values move. ---
$periods = @("2109","2110","2111","2112","2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# The "Valor Mora" (F) column keeps the same bag of values, but the one
# worker who had the smaller 33600 due amount now lines up with period
# 2211 (bottom row) instead of 2109 (top row).
$ws.Range("F16").Value = 48000
$ws.Range("F30").Value = 33600

# --- Nudge the company logo a bit to the left (its anchor cell offsets
# shift by 19pt / 241300 EMU while its size/top/height stay the same).
# Using absolute EMU-derived point values (off_x=680600 EMU, cx=975600 EMU)
# avoids drift from the Shape.Left getter using a different rounding model
# than the setter. ---
$logo = $ws.Shapes.Item(1)
$logo.Left = 53.59055118110236
$logo.Width = 76.81889763779527
